# modified test cases on overdue fix
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()

$wsSummary.Range("B2").Value = 822.33

# E2: style 8 -> 9, F2: style 9 -> 6 ; use format copy so the existing cellXfs
# entries are reused instead of creating new ones, then set the new values.
$wsSummary.Range("E3").Copy()
$wsSummary.Range("E2").PasteSpecial(-4122)
$wsSummary.Range("F2").Copy()
$wsSummary.Range("B2").PasteSpecial(-4122)
$wsSummary.Range("E2").Value = 4177.67
$wsSummary.Range("F2").Value = 851.5

# A3, E3, F3: style 10 -> 6 (keep values, fix format), then update the values
$wsSummary.Range("B3").Copy()
$wsSummary.Range("A3").PasteSpecial(-4122)
$wsSummary.Range("B3").Copy()
$wsSummary.Range("E3").PasteSpecial(-4122)
$wsSummary.Range("B3").Copy()
$wsSummary.Range("F3").PasteSpecial(-4122)

$wsSummary.Range("A3").Value = 111.25
$wsSummary.Range("B3").Value = 29.17
$wsSummary.Range("E3").Value = 82.08
$wsSummary.Range("F3").Value = 12.5

$wsSummary.Range("A7:XFD16").Select()

# ---------------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------------
$wsSched = $wb.Worksheets.Item("Repayment schedule")
$wsSched.Activate()

# Style fixes: many F/H (and a couple of G/K/P) cells move from the
# 2-decimal-place style (10) to the General style (6). Copy formats from a
# cell that already has style 6 so the existing cellXfs slot is reused.
$wsSched.Range("E3").Copy()
$wsSched.Range("F3,H3,F4,H4,F5,H5,F6,H6,F7,H7,F8,H8,F9,H9,K9,P9").PasteSpecial(-4122)
$wsSched.Range("E3").Copy()
$wsSched.Range("G8,G9").PasteSpecial(-4122)

# L3 moves from style 6 to style 10 (2-decimal-place)
$wsSched.Range("F3").Copy()
$wsSched.Range("L3").PasteSpecial(-4122)

# E4 moves from style 6 to style 12
$wsSched.Range("M3").Copy()
$wsSched.Range("E4").PasteSpecial(-4122)

# --- row 3 ---
$wsSched.Range("L3").Value = 839
$wsSched.Range("N3").Value = 839
$wsSched.Range("P3").Value = 12.5

# --- row 4 ---
$wsSched.Range("A4").Value = 2
$wsSched.Range("F4").Value = 0
$wsSched.Range("H4").Value = 12.5
$wsSched.Range("K4").Value = 12.5
$wsSched.Range("L4").Value = 12.5
$wsSched.Range("M4").Value = 0
$wsSched.Range("N4").Value = 0

# --- row 5 ---
$wsSched.Range("A5").Value = 3
$wsSched.Range("F5").Value = 839
$wsSched.Range("G5").Value = 3326.17
$wsSched.Range("H5").Value = 12.5

# --- row 6 ---
$wsSched.Range("A6").Value = 4
$wsSched.Range("F6").Value = 823.65
$wsSched.Range("G6").Value = 2502.52
$wsSched.Range("H6").Value = 27.85

# --- row 7 ---
$wsSched.Range("A7").Value = 5
$wsSched.Range("F7").Value = 830.65
$wsSched.Range("G7").Value = 1671.87
$wsSched.Range("H7").Value = 20.85

# --- row 8 ---
$wsSched.Range("A8").Value = 6
$wsSched.Range("F8").Value = 837.57
$wsSched.Range("G8").Value = 834.3
$wsSched.Range("H8").Value = 13.93

# --- row 9 ---
$wsSched.Range("A9").Value = 7
$wsSched.Range("F9").Value = 834.3
$wsSched.Range("G9").Value = 0
$wsSched.Range("H9").Value = 6.95
$wsSched.Range("K9").Value = 841.25
$wsSched.Range("P9").Value = 841.25

# row 12 (the stray helper cell F12) is removed entirely
$wsSched.Rows(12).Delete()

$wsSched.Range("L3").Select()

# ---------------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Activate()

$wsTrans.Range("E2").Copy()
$wsTrans.Range("F2,G2").PasteSpecial(-4122)

$wsTrans.Range("A2").Value = 64
$wsTrans.Range("E2").Value = 851.5
$wsTrans.Range("F2").Value = 822.33
$wsTrans.Range("G2").Value = 29.17
$wsTrans.Range("J2").Value = 4177.67

$wsTrans.Range("A3").Value = 62

$wsTrans.Range("A2:XFD4").Select()
